# Fruta / hortaliza, semanal
# Insert 3 new rows (new weekly price report, 2023-01-13) above the last
# "Femacal de La Calera" / Melon block, pushing the existing rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows before row 541; this shifts the former rows
# 541-549 down to 544-552, preserving their content & formatting.
$ws.Range("A541:R543").EntireRow.Insert()

# --- New row 541 : Tuna / Extra -------------------------------------------------
$ws.Cells.Item(541, 1).Value = 3
$ws.Cells.Item(541, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(541, 3).Value = "Coquimbo"
$ws.Cells.Item(541, 4).Value = 44939
$ws.Cells.Item(541, 5).Value = 5
$ws.Cells.Item(541, 6).Value = 100112027
$ws.Cells.Item(541, 7).Value = "Melón"
$ws.Cells.Item(541, 8).Value = "Tuna"
$ws.Cells.Item(541, 9).Value = "Extra"
$ws.Cells.Item(541, 10).Value = 1208
$ws.Cells.Item(541, 11).Value = 2000
$ws.Cells.Item(541, 12).Value = 2000
$ws.Cells.Item(541, 13).Value = 2000
$ws.Cells.Item(541, 14).Value = "$/unidad"
$ws.Cells.Item(541, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(541, 16).Value = 2000
$ws.Cells.Item(541, 17).Value = 1
$ws.Cells.Item(541, 18).Value = "Hortaliza"

# --- New row 542 : Tuna / Primera -----------------------------------------------
$ws.Cells.Item(542, 1).Value = 3
$ws.Cells.Item(542, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(542, 3).Value = "Coquimbo"
$ws.Cells.Item(542, 4).Value = 44939
$ws.Cells.Item(542, 5).Value = 5
$ws.Cells.Item(542, 6).Value = 100112027
$ws.Cells.Item(542, 7).Value = "Melón"
$ws.Cells.Item(542, 8).Value = "Tuna"
$ws.Cells.Item(542, 9).Value = "Primera"
$ws.Cells.Item(542, 10).Value = 1300
$ws.Cells.Item(542, 11).Value = 1500
$ws.Cells.Item(542, 12).Value = 1500
$ws.Cells.Item(542, 13).Value = 1500
$ws.Cells.Item(542, 14).Value = "$/unidad"
$ws.Cells.Item(542, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(542, 16).Value = 1500
$ws.Cells.Item(542, 17).Value = 1
$ws.Cells.Item(542, 18).Value = "Hortaliza"

# --- New row 543 : Tuna / Segunda -----------------------------------------------
$ws.Cells.Item(543, 1).Value = 3
$ws.Cells.Item(543, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(543, 3).Value = "Coquimbo"
$ws.Cells.Item(543, 4).Value = 44939
$ws.Cells.Item(543, 5).Value = 5
$ws.Cells.Item(543, 6).Value = 100112027
$ws.Cells.Item(543, 7).Value = "Melón"
$ws.Cells.Item(543, 8).Value = "Tuna"
$ws.Cells.Item(543, 9).Value = "Segunda"
$ws.Cells.Item(543, 10).Value = 1200
$ws.Cells.Item(543, 11).Value = 1000
$ws.Cells.Item(543, 12).Value = 1000
$ws.Cells.Item(543, 13).Value = 1000
$ws.Cells.Item(543, 14).Value = "$/unidad"
$ws.Cells.Item(543, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(543, 16).Value = 1000
$ws.Cells.Item(543, 17).Value = 1
$ws.Cells.Item(543, 18).Value = "Hortaliza"
